$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129, shifting existing rows 129-208 down to 130-209.
$ws.Rows("129:129").Insert()

# Populate the newly inserted row 129 with the new weekly record.
$ws.Range("A129").Value = 7
$ws.Range("B129").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C129").Value = "Ñuble"
$ws.Range("D129").Value = 45236
$ws.Range("E129").Value = 16
$ws.Range("F129").Value = "Fruta"
$ws.Range("G129").Value = 100108
$ws.Range("H129").Value = "Tropicales y subtropicales"
$ws.Range("I129").Value = 100108002
$ws.Range("J129").Value = "Mango"
$ws.Range("K129").Value = "Sin especificar"
$ws.Range("L129").Value = "Primera"
$ws.Range("M129").Value = 100
$ws.Range("N129").Value = 13000
$ws.Range("O129").Value = 13000
$ws.Range("P129").Value = 13000
$ws.Range("Q129").Value = "$/bandeja 4 kilos"
$ws.Range("R129").Value = "Brasil"
$ws.Range("S129").Value = 3250
$ws.Range("T129").Value = 4
